$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells: Wins / Losses / Ties in AD1:AF1
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Match the existing header formatting (bold, bordered, centered) by
# copying the format from the neighboring header cell (AC1) instead of
# rebuilding the font/border/alignment by hand - this keeps the new
# cells sharing the same style record as the rest of row 1.
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)  # xlPasteFormats

# Team record (Wins/Losses/Ties) repeated for every player row.
for ($r = 2; $r -le 57; $r++) {
    $ws.Cells.Item($r, 30).Value = 85   # AD - Wins
    $ws.Cells.Item($r, 31).Value = 77   # AE - Losses
    $ws.Cells.Item($r, 32).Value = 0    # AF - Ties
}
